$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the whole paragraph "Административный истец состоял на воинском
#    учете в Военном комиссариате." (entire <w:p> removed from the document).
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Административный истец состоял*на воинском учете в Военном комиссариате*") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2. "No" -> "№" typographic fix.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "04.07.2013 No 565", $true, $false, $false, $false, $false,
    $true, 1, $false, "04.07.2013 № 565", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Rewording: "...годности к военной службе. Возможности представить
#    имеющиеся у Административного истца медицинские документы не
#    предоставили." -> "...годности к военной службе. Административные
#    ответчики не дали возможности представить имеющиеся у Административного
#    истца медицинские документы."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Жалобы на состояние здоровья Административного истца не фиксировались, диагностические исследования не проводились, на обследование по имеющимся у Административного истца заболеваниям не направлялся несмотря на то, что без такого обследования и учета этих заболеваний невозможно вынесение заключения о категории годности к военной службе. Возможности представить имеющиеся у Административного истца медицинские документы не предоставили.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Жалобы на состояние здоровья Административного истца не фиксировались, диагностические исследования не проводились, на обследование по имеющимся у Административного истца заболеваниям не направлялся несмотря на то, что без такого обследования и учета этих заболеваний невозможно вынесение заключения о категории годности к военной службе. Административные ответчики не дали возможности представить имеющиеся у Административного истца медицинские документы.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 4. "Действия военного Административных ответчиков" -> "Действия
#    Административных ответчиков"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Действия военного Административных ответчиков незаконны",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Действия Административных ответчиков незаконны", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. "При доставлении в военкомат Административный истец написал..." ->
#    "Находясь в Военном комиссариате, Административный истец написал..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "При доставлении в военкомат Административный истец написал и требовал принять заявление о замене военной службы по призыву альтернативной гражданской службой. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Находясь в Военном комиссариате, Административный истец написал и требовал принять заявление о замене военной службы по призыву альтернативной гражданской службой.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 6. Insert two new paragraphs right after the "{%p endif %}" paragraph that
#    follows the alternative-civil-service-refusal passage (the one whose
#    previous sibling mentions "При доставлении в военкомат" / now "Находясь
#    в Военном комиссариате").
#
# NOTE: we deliberately ignore the object returned by InsertParagraphAfter()
# and instead re-enumerate $d.Paragraphs fresh (by index) after each
# insertion, since the returned handle does not reliably support further
# Range operations in this runtime.
# ---------------------------------------------------------------------------
$sourceIndex = 0
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Находясь в Военном комиссариате*") {
        $sourceIndex = $idx
    }
}

$anchorIndex = 0
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($idx -gt $sourceIndex -and $p.Range.Text -like "*endif*") {
        $anchorIndex = $idx
        break
    }
}

$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($idx -eq $anchorIndex) {
        $p.Range.InsertParagraphAfter() | Out-Null
    }
}

$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($idx -eq ($anchorIndex + 1)) {
        $p.Range.Text = "В силу статьи 28 Конституции России, статьи 18 Всеобщей декларации прав человека, статьи 18 Международного пакта о гражданских и политических правах, человек вправе придерживаться имеющихся убеждений, но и менять свою религию или убеждения."
    }
}

$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($idx -eq ($anchorIndex + 1)) {
        $p.Range.InsertParagraphAfter() | Out-Null
    }
}

$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($idx -eq ($anchorIndex + 2)) {
        $p.Range.Text = "В решении от 01.02.2013 по делу «Ким и др. против Республики Корея» (сообщение № 1786/2008 CCPR/C/106/D/1786/2008) Комитет ООН по правам человека приравнял право на отказ от военной службы к праву иметь убеждения, признав тем самым, что данное право не может быть ограничено государством (п. 1 ст. 18 Международного пакта о гражданских и политических правах)."
    }
}
